$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric value 10000 in G4 with the text "10 000",
# formatted the same as the other body cells (Aptos Narrow, 11pt, black),
# while keeping the existing cell number-format style (s="3").
$ws.Range("G4").Value = "10 000"
$c1 = $ws.Range("G4").Characters(1, 5)
$c1.Font.Name = "Aptos Narrow"
$c1.Font.Size = 11
$c1.Font.Color = 0
$c2 = $ws.Range("G4").Characters(6, 1)
$c2.Font.Name = "Aptos Narrow"
$c2.Font.Size = 11
$c2.Font.Color = 0

# Copy the fully-formatted G4 cell onto G6 so that both cells end up
# sharing the exact same shared-string entry, just like in the target file.
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
